$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "AA2" = 95
    "AJ2" = 320
    "AM2" = 580
    "AN2" = 120
    "G2" = 3.7
    "H2" = 2.32
    "P2" = 1.74
    "U2" = 1.98
    "W2" = 1.37
    "AC3" = 20
    "AG3" = 20
    "AJ3" = 180
    "AK3" = 70
    "G3" = 1.6
    "I3" = 9.800000000000001
    "N3" = 3.4
    "O3" = 1.31
    "S3" = 3
    "T3" = 2.02
    "X3" = 34
    "AD4" = 130
    "AH4" = 95
    "AJ4" = 10.5
    "AN4" = 6
    "G4" = 1.35
    "H4" = 8.199999999999999
    "I4" = 16
    "J4" = 4.9
    "L4" = 1.32
    "N4" = 3.75
    "O4" = 1.26
    "Q4" = 1.78
    "R4" = 1.42
    "T4" = 2.2
    "U4" = 1.63
    "V4" = 1.06
    "W4" = 3.75
    "Y4" = 100
    "AB5" = 29
    "AC5" = 12
    "AF5" = 18
    "AG5" = 40
    "AJ5" = 180
    "AK5" = 70
    "AN5" = 29
    "F5" = 1.53
    "G5" = 1.7
    "H5" = 5.3
    "I5" = 8
    "J5" = 3.85
    "K5" = 4.8
    "L5" = 1.4
    "M5" = 1.08
    "N5" = 3.15
    "O5" = 1.35
    "P5" = 1.82
    "Q5" = 1.94
    "R5" = 1.3
    "T5" = 2
    "U5" = 1.77
    "V5" = 1.15
    "W5" = 2.42
    "X5" = 90
    "AC6" = 10.5
    "AD6" = 13
    "AE6" = 28
    "F6" = 2.7
    "G6" = 2.98
    "J6" = 3.55
    "L6" = 1.34
    "N6" = 4.3
    "P6" = 2.12
    "R6" = 1.45
    "S6" = 2.86
    "U6" = 2.3
    "V6" = 1.59
    "AA7" = 10
    "AC7" = 21
    "AD7" = 13.5
    "AE7" = 15
    "AI7" = 160
    "AO7" = 3.45
    "F7" = 8.800000000000001
    "I7" = 1.25
    "J7" = 5.8
    "K7" = 9.6
    "L7" = 1.19
    "O7" = 1.13
    "P7" = 2.72
    "Q7" = 1.39
    "R7" = 1.78
    "S7" = 1.98
    "T7" = 2.02
    "U7" = 1.75
    "V7" = 4.8
    "X7" = 80
    "Y7" = 13.5
    "Z7" = 9.800000000000001
    "AE8" = 95
    "AK8" = 100
    "AL8" = 300
    "AM8" = 580
    "AO8" = 38
    "J8" = 3.05
    "N8" = 3
    "Q8" = 2.12
    "AE9" = 330
    "AF9" = 11.5
    "AH9" = 19.5
    "AI9" = 330
    "AJ9" = 21
    "AO9" = 80
    "J9" = 3.75
    "N9" = 3.65
    "P9" = 1.91
    "Q9" = 1.98
    "R9" = 1.35
    "Y9" = 17
    "AA10" = 60
    "AB10" = 14.5
    "AD10" = 14
    "AE10" = 34
    "AF10" = 17.5
    "AH10" = 17
    "AI10" = 38
    "AJ10" = 34
    "AK10" = 26
    "AM10" = 330
    "AN10" = 15.5
    "F10" = 2.46
    "G10" = 2.52
    "H10" = 2.92
    "J10" = 3.7
    "K10" = 3.9
    "L10" = 1.3
    "M10" = 1.04
    "N10" = 5.1
    "O10" = 1.22
    "P10" = 2.32
    "Q10" = 1.7
    "R10" = 1.5
    "S10" = 2.78
    "T10" = 1.6
    "U10" = 2.52
    "V10" = 1.48
    "W10" = 1.66
    "X10" = 22
    "Y10" = 16.5
    "Z10" = 22
    "AA11" = 44
    "AB11" = 23
    "AC11" = 10.5
    "AD11" = 10.5
    "AE11" = 17
    "AF11" = 34
    "AG11" = 16.5
    "AH11" = 15.5
    "AI11" = 26
    "AJ11" = 95
    "AK11" = 120
    "AL11" = 160
    "AN11" = 27
    "F11" = 3.9
    "G11" = 4
    "H11" = 1.92
    "I11" = 1.94
    "J11" = 4.3
    "K11" = 4.5
    "L11" = 1.28
    "N11" = 5.7
    "O11" = 1.19
    "P11" = 2.68
    "Q11" = 1.55
    "R11" = 1.66
    "S11" = 2.38
    "T11" = 1.57
    "U11" = 2.62
    "V11" = 2.06
    "W11" = 1.33
    "Y11" = 14
    "Z11" = 16
    "AA12" = 44
    "AM12" = 580
    "AN12" = 34
    "G12" = 2.96
    "L12" = 1.41
    "P12" = 1.83
    "U12" = 2.2
    "W12" = 1.51
    "Z12" = 23
    "AB13" = 12
    "AE13" = 130
    "AF13" = 13.5
    "AJ13" = 21
    "AL13" = 32
    "AM13" = 580
    "G13" = 1.92
    "O13" = 1.23
    "V13" = 1.25
    "Z13" = 85
    "AB14" = 8.800000000000001
    "AC14" = 8.800000000000001
    "AD14" = 22
    "AE14" = 1000
    "AH14" = 26
    "AI14" = 1000
    "AL14" = 50
    "AN14" = 20
    "F14" = 1.91
    "G14" = 1.97
    "I14" = 4.7
    "J14" = 3.6
    "K14" = 3.95
    "W14" = 2.02
    "X14" = 13.5
    "Z14" = 36
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
